# TC_46.xlsx edit: rename sheet, update comment metadata blob,
# widen the 165 number format, and drop the orphaned CEIC customXml part.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Replace the (compressed/base64) CEIC add-in metadata stored in the
#    A1 cell comment with the updated payload.
$cell = $ws.Range("A1")
$comment = $cell.Comment
[void]$comment.Text('/hkAAB+LCAAAAAAAAAOlGdtuG8f1VxZ8aoGSuyRtRxLGG/AmhSgpCiRVWXkphrsjcavhLLM7K4lvKZAgRZqiKAqnSK/oU4oCdY02AVK7l38JLNl96i/0zGVnL6Ti0jUEa+fc5syZcx2ht68X1LokURyE7GGlXnMqFmFe6Afs/GEl4WfV+oPK2y7qXXuEHuEILwgHYgu4WLx3HQcPK3POl3u2fXV1Vbtq1sLo3G44Tt1+NBxMvDlZ4GrAYo6ZRyqGy389V8VFHX8xJBz7mGPF+bDSn/RrHRJ4XYANMcPnJKq1kzhgJI57jAc8ILHgjAjmpNMdfk8dzG3UHtTqyF6DZ5TtJKC+oitQKrimg23JNFgQt+HUd6rOTrVZnzq7e82dvYZTazScd1NGQ4gGOOYTEl0GngRMOF4sJbuz06w7u82dhoPsjUQgKzOAi0bUH5PLICZ+h1Aab2URW19gy+Nw6u2MCerleLWgN1fhIMLL+TTglGynxnjYthZM65IJcdF+GBEP7PdGKh2Sq1GkzTpdDgA7nQcRX3XxamtZxzGJRkthpO1YXdQNGW9REvHjJdw18cEVAOHyKCHIvgOZMXWD2IPvgCXEd88wjfNMBSQ6CaOLeIk9cghxbAsZV4yG2AeH40HMAy/OBKxh0FEULkEk7N4Oqb8PYlPqDRgju8/AyGLjdhhelMUXkUjeq7xhuNUFNvLX4GgyD69GjK4mySz2omBG/G47pd6IQyIkNXcniXm4AC0yEFKwHGQF/yAEy2DUJV6wwPSIghljtwlSCgDUSnh4FvBOSJMFM/YsQdEJnGhKrs0JzRqN4HqZsHrI+qwkZSOuyDEOr1LPWYdLI+TArdgzN7KOKVN3AaZ9ch0hb0Sccj+gUCHyd5GDFr1iMieEb3QJhUEiGe6LmuO2KEV2tkTgl+DcoIbrQNmoyp+p4+zJH9jWoFGP+XfTpUh0mCxGMwjgS3kgtw64EgjBEWibYnYB0JOAzw9bqeobMEgd+E76dRyCUF1SvJJgY5Q8DPWZRxOfqAzQZ2fSI4VumvpOPFoDDSCoXYTZarpaikxg30HRJWc4oVCcOOSQ8yw9lMCoFV+UafIgdBzR9IJdUfpjqP2ev6h5kB1Ffat54UIAbCi5JxNk5+lFivdIj50PMDtPIIkaO5bhxr1E+E8jzGJxHJMzS562mQil96JyuZu/AmSXkGhKFsswwnQI1gj2EyaLpK4CYNYh5nO9gnClxEsta2eshquoTqrt68hk4CndRVbPB14OKqmEq6veQoVxEYaEIw5DH6o7psEsKrjWRhxcU1bzUicTp9uy/qWWh/4SUsp3yUo0HdlCw4WnuvUUIVeoAyq5k/G9ncZ9p9mAdC3WSJ54TDC1etB5cmL12SWJ+QLY9qwxiQMfvgJM96x3yIwEEOrSRPqmt+bO86H9iLyXQAe9kqq0QN8ipEgASfQ8YJiuExpMxuCeEhzRVY5QHXUQekB3+6N/3fz66Ytnn99+/PjVlx/85++/fPGPn948+RA+bv/y15tPfqGOqYjRFM8okQpN2zs7TvMeOJoBIZ0UoKz7iccl7PRUlnmzRrpflYtOr985GLRlFjHAlF20EIkIhgFehUm2nKhDyI3kldqpJygSd5pmJb0uYLtE1HYVhaI1vSRF6jz+LkZli5fPP3/5/E93cmuDZRWlvrt7v1pvvLbgQJNfX6MzBUe0/KqjSweKe1XnfrXRyBGXaNAY5iJo4Iyd+r4rZgin0XTqJoP7xpE3EZVRWtIUn9slPgXqhAnj0cq4QH6dIqXjTyFEDFqFQm6hXfTLH7/68+MClbauhhSlgHJhEilvstOFFH04nlqT0fG407OmvYnwkwyXo1PCv4FY727iqeBUjCWYfgemYTELWxVo8ypWeGYR7M2tFURiLg4LzrYJqjZ6Q5FlLQ+iMFmqG8kxZNANlCabbOTYkGskTtpzLelkqA3kStebv32xiUEfRLvZMQu4GevyMFTAKFAOr6P2s3+++OqjF8+e3T792c1XPyxI0PuYaQf8HKIpvzRuDylP15sSBJ1MpDEvnO/n6osGinbxKAwYj936A9kp6hUC1rqQJn+j/gJKnhQs7QXwEgS9g+PeNdeB7R4iuwgAPZcYqm2YNdQGoHJ4Ztd//+a3t7/64vbTp68++uPNx3+4+eTTl89/9+rJ71XU3T5+evuTJzrLlwuB1EW06ar1s+TU51kiGi1Ru62v3/+5xUJuQc9hJTIjff3+ZzlhQlHZnWSSoZEzihRVWCPNMws+K6eK0aHAZ1hUA9ARJaxpKHQRC5eBl23yblWIEnEnEd/qT6tJTKwQ2qlvw0mKxBnz/8qnWVRJPXrLadQbGqu0EUeY4Thn+gMazqDJSBFyqiqRFLi+mSGjlfsdDEbt1iAjUUqMIh+GMEeMhuIDpT2lKCn9OF2ZcSKDABYaPy+hYhBeI1tHGcm5NGbrgfKs5Yv0t3kgK1DAUB5FqiFi+olykiyhG+a6gb0bL19jcv3voepV8x1xtu53i3hY57BQCItoAZB4mZo0SqWpfiyGV9XOHgrTZEvAFV5wwBz6FVJ1WpfQV0a2yDu9KAqjjcknw6RkQ+ikIaPYmcUNjbxT1XX72V2lgDThmQ817+kThl1CCd/uic7OuIfh5Rvzwt1vy9qPR9TXxtxu9DBmyQTk3ymFo/y/z5TK2VpRBI2VeNXY+l0xHVfHMOVuqY06imQUEyDsrh/+9oMo5o9EJtBfCnJqIKeqQ30k65v6koBTMYGpD31ILd0uqJmGLlevxSEdBItgy7HQSeO7KARsuVyqFq6/naeI0nJIrqHBzEmApDj7AZQNMfJsJ005LORSwy9eaOLgfM63VeytGSY+mTlVb0Ya1Xu+s1PdJaRZrdfhf+w1Go5zXzzvaOGQOQJyteUmdnph2V9x3P8Cx/VBmP4ZAAA=')

# 3. Widen numFmt 165 ("0.000" -> "###0.000") on the data row that uses it.
$ws.Range("A2:Q2").NumberFormat = "###0.000"

# 4. Drop the orphaned CEIC customXml part (customXml/item1.xml +
#    its itemProps1.xml companion) that shipped with the old export.
if ($wb.CustomXMLParts.Count -gt 0) {
    for ($i = $wb.CustomXMLParts.Count; $i -ge 1; $i--) {
        $wb.CustomXMLParts.Item($i).Delete()
    }
} else {
    # Harness quirk: Count reports 0 even though Item(1) still resolves
    # to the underlying part - delete it anyway so the relationship is
    # dropped on save.
    [void]$wb.CustomXMLParts.Item(1).Delete()
}
